$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formats from row 25 (a "normal" data row) for the two new non-final
# rows, and from row 26 (the previous "last row", using the heavier
# last-row style) for the new final row, so the new rows inherit the same
# style indices used throughout the table.
$ws.Range("A25:E25").Copy()
$ws.Range("A27:E28").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A26:E26").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122) # xlPasteFormats

# Row 27 - RCC090
$ws.Range("A27").Value = "RCC090"
$ws.Range("B27").Value = "OPQA-1467||OPQA-1518"
$ws.Range("C27").Value = "Verify that Group owner user is able to add cover photo to a group || Verify that cover photo is displayed properly from groups tab for the newly created group."
$ws.Range("D27").Value = "Y"

# Row 28 - RCC091
$ws.Range("A28").Value = "RCC091"
$ws.Range("B28").Value = "OPQA-1581"
$ws.Range("C28").Value = "Verify that group owner is able to modify the cover photo of the group through Edit option from group details page."
$ws.Range("D28").Value = "Y"

# Row 29 - RCC092
$ws.Range("A29").Value = "RCC092"
$ws.Range("B29").Value = "OPQA-1584"
$ws.Range("C29").Value = "Verify that group owner sees option to add an image if cover photo not exists for the group from group details page."
$ws.Range("D29").Value = "Y"

$ws.Range("A27:E27").RowHeight = 30

$ws.Range("C34").Select()
